$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Reuse the exact header formatting (bold, border, center/top alignment)
# already used by the other header cells (e.g. AA1) instead of building a
# brand-new style.
$ws.Range("AA1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-43
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 103   # AD
    $ws.Cells.Item($r, 31).Value = 59    # AE
    $ws.Cells.Item($r, 32).Value = 0     # AF
}
